# Upgrade navigation menu accessibility
# Append a new data row (row 84) to each of the 4 worksheets, mirroring
# the existing row layout (time, 总长, ID, 实际长度, 和校验, and the
# corresponding *_DEC numeric columns).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newRow = 84

# --- Sheet 1: FE_LFT_#1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item($newRow, 1).Value = 45870.49162037037
$ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat
$ws.Cells.Item($newRow, 2).Value = "0x01,0x7c"
$ws.Cells.Item($newRow, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item($newRow, 4).Value = "0x01,0x18"
$ws.Cells.Item($newRow, 5).Value = "0xf"
$ws.Cells.Item($newRow, 6).Value = 380
$ws.Cells.Item($newRow, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item($newRow, 8).Value = 280
$ws.Cells.Item($newRow, 9).Value = 15

# --- Sheet 2: FE_LFT_#2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item($newRow, 1).Value = 45870.49162037037
$ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat
$ws.Cells.Item($newRow, 2).Value = "0x01,0x90"
$ws.Cells.Item($newRow, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item($newRow, 4).Value = "0x01,0x24"
$ws.Cells.Item($newRow, 5).Value = "0xe"
$ws.Cells.Item($newRow, 6).Value = 400
$ws.Cells.Item($newRow, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item($newRow, 8).Value = 292
$ws.Cells.Item($newRow, 9).Value = 14

# --- Sheet 3: FE_PLT_#1 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item($newRow, 1).Value = 45870.49162037037
$ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat
$ws.Cells.Item($newRow, 2).Value = "0x00,0x6e"
$ws.Cells.Item($newRow, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item($newRow, 4).Value = "0x00,0x5F"
$ws.Cells.Item($newRow, 5).Value = "0x3"
$ws.Cells.Item($newRow, 6).Value = 110
$ws.Cells.Item($newRow, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($newRow, 8).Value = 95
$ws.Cells.Item($newRow, 9).Value = 3

# --- Sheet 4: FE_PLT_#2 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item($newRow, 1).Value = 45870.49162037037
$ws.Cells.Item($newRow, 1).NumberFormat = $dateFormat
$ws.Cells.Item($newRow, 2).Value = "0x00,0x6e"
$ws.Cells.Item($newRow, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item($newRow, 4).Value = "0x00,0x5D"
$ws.Cells.Item($newRow, 5).Value = "0x3"
$ws.Cells.Item($newRow, 6).Value = 110
$ws.Cells.Item($newRow, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item($newRow, 8).Value = 93
$ws.Cells.Item($newRow, 9).Value = 3
